# Auto-generated Excel COM-interop script applying crypto price/volume updates
# (commit: "Updated cryptos list ... with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.921.77'
$ws.Range('E2').Value = '  +0.00%  '
$ws.Range('D3').Value = '1.629.88'
$ws.Range('E3').Value = '  -0.03%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '''211.79'
$ws.Range('E5').Value = '  -0.08%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  -0.46%  '
$ws.Range('E9').Value = '  +0.34%  '
$ws.Range('E10').Value = '  -0.93%  '
$ws.Range('D11').Value = '''0.0880'
$ws.Range('E11').Value = '  -0.02%  '
$ws.Range('D12').Value = '1.861.89'
$ws.Range('E12').Value = '  +0.02%  '
$ws.Range('D13').Value = '1.616.48'
$ws.Range('E13').Value = '  -0.75%  '
$ws.Range('E14').Value = '  -0.92%  '
$ws.Range('E15').Value = '  -1.05%  '
$ws.Range('D16').Value = '''64.87'
$ws.Range('D17').Value = '27.930.30'
$ws.Range('E17').Value = '  +0.04%  '
$ws.Range('D18').Value = '''228.11'
$ws.Range('E18').Value = '  -0.88%  '
$ws.Range('E20').Value = '  -0.66%  '
$ws.Range('E21').Value = '  +0.01%  '
$ws.Range('E22').Value = '  +0.03%  '
$ws.Range('E23').Value = '  -2.82%  '
$ws.Range('E24').Value = '  +1.93%  '
$ws.Range('D25').Value = '''154.65'
$ws.Range('E25').Value = '  -0.17%  '
$ws.Range('E26').Value = '  -0.26%  '
$ws.Range('E27').Value = '  -0.36%  '
$ws.Range('E28').Value = '  -0.05%  '
$ws.Range('D29').Value = '''15.40'
$ws.Range('E29').Value = '  -0.99%  '
$ws.Range('E30').Value = '  -0.25%  '
$ws.Range('D31').Value = '''0.0481'
$ws.Range('E31').Value = '  -0.25%  '
$ws.Range('E32').Value = '  +0.30%  '
$ws.Range('D33').Value = '1.419.97'
$ws.Range('E33').Value = '  +1.46%  '
$ws.Range('E34').Value = '  +1.09%  '
$ws.Range('D35').Value = '''1.63'
$ws.Range('E35').Value = '  +3.31%  '
$ws.Range('E36').Value = '  -1.56%  '
$ws.Range('E37').Value = '  -1.04%  '
$ws.Range('E38').Value = '  -0.70%  '
$ws.Range('D39').Value = '''0.554'
$ws.Range('E39').Value = '  -0.16%  '
$ws.Range('D40').Value = '''0.852'
$ws.Range('E40').Value = '  -1.35%  '
$ws.Range('E41').Value = '  -2.74%  '
$ws.Range('D42').Value = '''65.75'
$ws.Range('E42').Value = '  -1.00%  '
$ws.Range('E43').Value = '  -1.27%  '
$ws.Range('E44').Value = '  -0.77%  '
$ws.Range('D45').Value = '1.770.82'
$ws.Range('E45').Value = '  -0.08%  '
$ws.Range('E46').Value = '  -3.67%  '
$ws.Range('D47').Value = '''88.76'
$ws.Range('E47').Value = '  +0.57%  '
$ws.Range('E48').Value = '  +0.67%  '
$ws.Range('E49').Value = '  -0.44%  '
$ws.Range('D50').Value = '''7.58'
$ws.Range('E50').Value = '  +0.49%  '
$ws.Range('E51').Value = '  -0.01%  '
